$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell reference -> new literal text value.
# Cells are written as plain text (matching the original inlineStr string
# cells) by temporarily forcing a text number format before assigning the
# value, then restoring the cell style so no visible formatting changes.
$updates = [ordered]@{
    'D2' = '306.04'
    'E2' = '1.22%'
    'D3' = '36.32'
    'E3' = '-0.99%'
    'D4' = '5.058'
    'E4' = '1.57%'
    'D5' = '0.07937'
    'E5' = '3.02%'
    'D6' = '2.238'
    'E6' = '6.30%'
    'D7' = '8.007'
    'E7' = '0.76%'
    'B8' = 'GateToken'
    'C8' = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
    'D8' = '4.150'
    'E8' = '3.11%'
    'B9' = 'MXToken'
    'C9' = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
    'D9' = '0.9276'
    'E9' = '1.61%'
    'B10' = 'LiechtensteinCryptoassetsExchange'
    'C10' = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
    'D10' = '0.09835'
    'E10' = '1.66%'
    'B11' = 'WazirX'
    'C11' = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
    'D11' = '0.1879'
    'E11' = '1.66%'
    'B12' = 'MandalaExchangeToken'
    'C12' = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
    'D12' = '0.09096'
    'E12' = '6.22%'
    'B13' = 'BitrueCoin'
    'C13' = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
    'D13' = '0.03713'
    'E13' = '4.51%'
    'B14' = 'BitMartToken'
    'C14' = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
    'D14' = '0.09914'
    'E14' = '-0.40%'
    'B15' = 'BitForexToken'
    'C15' = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
    'D15' = '0.001430'
    'E15' = '-2.46%'
    'B16' = 'TigerCash'
    'C16' = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
    'D16' = '0.005627'
    'E16' = '-0.88%'
    'B17' = 'LEO'
    'C17' = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
    'D17' = '3.463'
    'E17' = '-0.07%'
    'D18' = '2.632'
    'E18' = '18.59%'
    'E19' = '-0.47%'
    'D20' = '0.1319'
    'E20' = '-0.71%'
    'D21' = '5.094'
    'E21' = '7.16%'
    'D22' = '0.2246'
    'E22' = '2.27%'
    'E23' = '-0.96%'
    'D24' = '0.001239'
    'E24' = '0.87%'
    'D25' = '0.004779'
    'E25' = '-6.36%'
    'D26' = '0.0001299'
    'E26' = '-7.07%'
    'D39' = '0.01920'
    'E39' = '8.98%'
    'D40' = '0.04927'
    'E40' = '6.73%'
    'D41' = '0.007793'
    'E41' = '4.49%'
    'D42' = '0.1393'
    'E42' = '0.26%'
    'D43' = '0.007796'
    'E43' = '0.93%'
    'D44' = '0.002124'
    'E44' = '-1.51%'
    'D45' = '0.01141'
    'E45' = '10.28%'
    'D46' = '0.00006246'
    'E46' = '-0.83%'
    'E47' = '0.08%'
    'D48' = '52.29'
    'E48' = '40.73%'
    'E49' = '-9.94%'
    'E50' = '0.08%'
    'E51' = '0.08%'
}

foreach ($cellRef in $updates.Keys) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = '@'
    $cell.Value = $updates[$cellRef]
    $cell.Style = 'Normal'
}

Write-Host "Applied $($updates.Count) cell updates"
